# Applies the Lamia_Profits scheduled-runner price/profit update.
# Generated from the canonical OOXML diff: updates currentAveragePrice /
# NQ/HQ price & profit columns (H:N) on specific rows across all 8 sheets.

$wb = $excel.ActiveWorkbook

# Sheet ALC row 12
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(12, 8).Value = 373.125
$ws.Cells.Item(12, 10).Value = 409.75
$ws.Cells.Item(12, 12).Value = 409.75
$ws.Cells.Item(12, 14).Value = -749.75

# Sheet ALC row 64
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(64, 8).Value = 6497.5557
$ws.Cells.Item(64, 9).Value = 4210.7144
$ws.Cells.Item(64, 11).Value = 4210.7144
$ws.Cells.Item(64, 13).Value = -3962.7144

# Sheet ALC row 67
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(67, 8).Value = 6497.5557
$ws.Cells.Item(67, 9).Value = 4210.7144
$ws.Cells.Item(67, 11).Value = 4210.7144
$ws.Cells.Item(67, 13).Value = -3352.7144

# Sheet ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 11238726
$ws.Cells.Item(137, 9).Value = 37038896
$ws.Cells.Item(137, 11).Value = 111116688
$ws.Cells.Item(137, 13).Value = -111114138

# Sheet ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138, 8).Value = 3669.8552
$ws.Cells.Item(138, 9).Value = 1940.5
$ws.Cells.Item(138, 10).Value = 4131.0166
$ws.Cells.Item(138, 11).Value = 5821.5
$ws.Cells.Item(138, 12).Value = 12393.0498
$ws.Cells.Item(138, 13).Value = -681.5
$ws.Cells.Item(138, 14).Value = -22673.0498

# Sheet ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 5407.0835
$ws.Cells.Item(2, 9).Value = 1534.4166
$ws.Cells.Item(2, 11).Value = 1534.4166
$ws.Cells.Item(2, 13).Value = -1421.4166

# Sheet ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 5862.3833
$ws.Cells.Item(32, 9).Value = 5862.3833
$ws.Cells.Item(32, 11).Value = 5862.3833
$ws.Cells.Item(32, 13).Value = -5575.3833

# Sheet ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 3277.2666
$ws.Cells.Item(45, 9).Value = 3262.5
$ws.Cells.Item(45, 11).Value = 3262.5
$ws.Cells.Item(45, 13).Value = -2885.5

# Sheet ARM row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(110, 8).Value = 3118.1304
$ws.Cells.Item(110, 9).Value = 2682.25
$ws.Cells.Item(110, 11).Value = 2682.25
$ws.Cells.Item(110, 13).Value = -637.25

# Sheet ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(116, 8).Value = 5407.0835
$ws.Cells.Item(116, 9).Value = 1534.4166
$ws.Cells.Item(116, 11).Value = 1534.4166
$ws.Cells.Item(116, 13).Value = 759.5834

# Sheet ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(122, 8).Value = 3016.6943
$ws.Cells.Item(122, 9).Value = 2460.7083
$ws.Cells.Item(122, 11).Value = 7382.124899999999
$ws.Cells.Item(122, 13).Value = -4932.124899999999

# Sheet BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 5407.0835
$ws.Cells.Item(3, 9).Value = 1534.4166
$ws.Cells.Item(3, 11).Value = 1534.4166
$ws.Cells.Item(3, 13).Value = -1420.4166

# Sheet BSM row 59
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(59, 8).Value = 89998.5
$ws.Cells.Item(59, 10).Value = 89998.5
$ws.Cells.Item(59, 12).Value = 89998.5
$ws.Cells.Item(59, 14).Value = -91692.5

# Sheet BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 3653.1714
$ws.Cells.Item(86, 10).Value = 5790.875
$ws.Cells.Item(86, 12).Value = 5790.875
$ws.Cells.Item(86, 14).Value = -8036.875

# Sheet BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(89, 8).Value = 3653.1714
$ws.Cells.Item(89, 10).Value = 5790.875
$ws.Cells.Item(89, 12).Value = 28954.375
$ws.Cells.Item(89, 14).Value = -40186.375

# Sheet BSM row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 1407.8889
$ws.Cells.Item(94, 9).Value = 1495.875
$ws.Cells.Item(94, 11).Value = 1495.875
$ws.Cells.Item(94, 13).Value = -1044.875

# Sheet BSM row 101
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(101, 8).Value = 75000
$ws.Cells.Item(101, 10).Value = 75000
$ws.Cells.Item(101, 12).Value = 75000
$ws.Cells.Item(101, 14).Value = -81490

# Sheet BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 1931.8667
$ws.Cells.Item(107, 9).Value = 2303.2222
$ws.Cells.Item(107, 10).Value = 1374.8334
$ws.Cells.Item(107, 11).Value = 2303.2222
$ws.Cells.Item(107, 12).Value = 1374.8334
$ws.Cells.Item(107, 13).Value = -383.2222000000002
$ws.Cells.Item(107, 14).Value = -5214.8334

# Sheet CRP row 51
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(51, 8).Value = 39999.715

# Sheet CRP row 61
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(61, 8).Value = 39999.715

# Sheet CRP row 68
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(68, 8).Value = 57776.5
$ws.Cells.Item(68, 10).Value = 57776.5
$ws.Cells.Item(68, 12).Value = 57776.5
$ws.Cells.Item(68, 14).Value = -59274.5

# Sheet CRP row 71
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(71, 8).Value = 57776.5
$ws.Cells.Item(71, 10).Value = 57776.5
$ws.Cells.Item(71, 12).Value = 173329.5
$ws.Cells.Item(71, 14).Value = -180817.5

# Sheet CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 4012.48
$ws.Cells.Item(132, 9).Value = 3554.0833
$ws.Cells.Item(132, 10).Value = 15014
$ws.Cells.Item(132, 11).Value = 10662.2499
$ws.Cells.Item(132, 12).Value = 45042
$ws.Cells.Item(132, 13).Value = -8132.249899999999
$ws.Cells.Item(132, 14).Value = -50102

# Sheet CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(134, 8).Value = 3396.55
$ws.Cells.Item(134, 9).Value = 2407.4707
$ws.Cells.Item(134, 11).Value = 7222.4121
$ws.Cells.Item(134, 13).Value = -4687.4121

# Sheet CUL row 61
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(61, 8).Value = 1537.25
$ws.Cells.Item(61, 9).Value = 199.33333
$ws.Cells.Item(61, 10).Value = 2340
$ws.Cells.Item(61, 11).Value = 597.99999
$ws.Cells.Item(61, 12).Value = 7020
$ws.Cells.Item(61, 13).Value = -382.99999
$ws.Cells.Item(61, 14).Value = -7450

# Sheet CUL row 75
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(75, 8).Value = 142861700
$ws.Cells.Item(75, 9).Value = 500000500
$ws.Cells.Item(75, 11).Value = 1500001500
$ws.Cells.Item(75, 13).Value = -1500000502

# Sheet CUL row 78
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(78, 8).Value = 142861700
$ws.Cells.Item(78, 9).Value = 500000500
$ws.Cells.Item(78, 11).Value = 4500004500
$ws.Cells.Item(78, 13).Value = -4499999508

# Sheet CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 5154155.5
$ws.Cells.Item(131, 9).Value = 20834288
$ws.Cells.Item(131, 10).Value = 3474141.5
$ws.Cells.Item(131, 11).Value = 62502864
$ws.Cells.Item(131, 12).Value = 10422424.5
$ws.Cells.Item(131, 13).Value = -62497824
$ws.Cells.Item(131, 14).Value = -10432504.5

# Sheet CUL row 137
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(137, 8).Value = 126712.625
$ws.Cells.Item(137, 9).Value = 1957.1428
$ws.Cells.Item(137, 10).Value = 1000001
$ws.Cells.Item(137, 11).Value = 5871.428400000001
$ws.Cells.Item(137, 12).Value = 3000003
$ws.Cells.Item(137, 13).Value = -771.4284000000007
$ws.Cells.Item(137, 14).Value = -3010203

# Sheet GSM row 97
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 1782.5264
$ws.Cells.Item(97, 10).Value = 8165.3335
$ws.Cells.Item(97, 12).Value = 8165.3335
$ws.Cells.Item(97, 14).Value = -9157.333500000001

# Sheet GSM row 106
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(106, 8).Value = 71250
$ws.Cells.Item(106, 10).Value = 71250
$ws.Cells.Item(106, 12).Value = 71250
$ws.Cells.Item(106, 14).Value = -73774

# Sheet LTW row 16
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 1469.6471
$ws.Cells.Item(16, 9).Value = 691.53845
$ws.Cells.Item(16, 10).Value = 3998.5
$ws.Cells.Item(16, 11).Value = 691.53845
$ws.Cells.Item(16, 12).Value = 3998.5
$ws.Cells.Item(16, 13).Value = -521.53845
$ws.Cells.Item(16, 14).Value = -4338.5

# Sheet LTW row 112
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(112, 8).Value = 0
$ws.Cells.Item(112, 10).Value = 0
$ws.Cells.Item(112, 12).Value = 0
$ws.Cells.Item(112, 14).ClearContents()

# Sheet LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 6909.4165
$ws.Cells.Item(132, 9).Value = 6909.4165
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 11).Value = 20728.2495
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 13).ClearContents()
$ws.Cells.Item(132, 14).Value = -18198.2495

# Sheet WVR row 82
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(82, 8).Value = 49997.5
$ws.Cells.Item(82, 9).Value = 49997.5
$ws.Cells.Item(82, 10).Value = 0
$ws.Cells.Item(82, 11).Value = 49997.5
$ws.Cells.Item(82, 12).Value = 0
$ws.Cells.Item(82, 13).ClearContents()
$ws.Cells.Item(82, 14).Value = -49614.5

# Sheet WVR row 85
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(85, 8).Value = 49997.5
$ws.Cells.Item(85, 9).Value = 49997.5
$ws.Cells.Item(85, 10).Value = 0
$ws.Cells.Item(85, 11).Value = 49997.5
$ws.Cells.Item(85, 12).Value = 0
$ws.Cells.Item(85, 13).ClearContents()
$ws.Cells.Item(85, 14).Value = -48671.5

# Sheet WVR row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 2732.1
$ws.Cells.Item(126, 9).Value = 1743.1305
$ws.Cells.Item(126, 10).Value = 5981.5713
$ws.Cells.Item(126, 11).Value = 5229.3915
$ws.Cells.Item(126, 12).Value = 17944.7139
$ws.Cells.Item(126, 13).Value = -2759.3915
$ws.Cells.Item(126, 14).Value = -22884.7139

# Sheet WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 2116.6775
$ws.Cells.Item(132, 9).Value = 1576.2307
$ws.Cells.Item(132, 10).Value = 4927
$ws.Cells.Item(132, 11).Value = 4728.6921
$ws.Cells.Item(132, 12).Value = 14781
$ws.Cells.Item(132, 13).Value = -2198.6921
$ws.Cells.Item(132, 14).Value = -19841

